$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume(1h) figures from the latest symbol-list refresh.
# Cells are stored as text (matching the source data format), so force a Text
# number format before writing to prevent Excel from auto-converting the
# numeric-looking / percentage-looking strings into actual numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.58%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.04%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08162"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.052"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.944"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.128"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.38%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.852"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "7.36%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9252"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.12%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1099"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "13.66%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1910"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09172"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.27%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.48%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09906"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.04%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001436"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.64%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005690"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.74%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.471"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.03%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3393"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.22%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.86%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.094"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2213"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.45%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04544"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.51%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.62%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004781"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.54%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.66%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004450"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01961"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.13%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.87%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007574"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.64%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009976"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "29.08%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.44%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002198"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01161"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.90%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006564"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.77%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.18%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "60.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.32%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001501"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-20.92%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.18%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.18%"
